$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly price block (row 176),
# pushing the existing data (old rows 176-296) down to 178-298.
$ws.Range("A176:R177").Insert()

# New week's data: row 176 = "Primera" quality, row 177 = "Segunda" quality.
$ws.Range("A176").Value = 1
$ws.Range("B176").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C176").Value = "Arica y Parinacota"
$ws.Range("D176").Value = 44673
$ws.Range("E176").Value = 15
$ws.Range("F176").Value = 100114014
$ws.Range("G176").Value = "Betarraga"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 1000
$ws.Range("K176").Value = 450
$ws.Range("L176").Value = 500
$ws.Range("M176").Value = 475
$ws.Range("N176").Value = "`$/paquete 4 unidades"
$ws.Range("O176").Value = "Región de Arica y Parinacota"
$ws.Range("P176").Value = 119
$ws.Range("Q176").Value = 4
$ws.Range("R176").Value = "Hortaliza"

$ws.Range("A177").Value = 1
$ws.Range("B177").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C177").Value = "Arica y Parinacota"
$ws.Range("D177").Value = 44673
$ws.Range("E177").Value = 15
$ws.Range("F177").Value = 100114014
$ws.Range("G177").Value = "Betarraga"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Segunda"
$ws.Range("J177").Value = 1200
$ws.Range("K177").Value = 450
$ws.Range("L177").Value = 500
$ws.Range("M177").Value = 475
$ws.Range("N177").Value = "`$/paquete 5 unidades"
$ws.Range("O177").Value = "Región de Arica y Parinacota"
$ws.Range("P177").Value = 95
$ws.Range("Q177").Value = 5
$ws.Range("R177").Value = "Hortaliza"
